$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values (column D) are stored as text in the source data (e.g. "28.718.07",
# "0.4620" with a significant trailing zero), so force text format before assigning
# to avoid Excel auto-converting numeric-looking strings to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.718.07"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.872.53"
$ws.Range("E3").Value = "  +2.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.76"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4620"
$ws.Range("E7").Value = "  -0.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3881"
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07878"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9764"
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.87"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.860.37"
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.015"
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.703"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06965"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.46"
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.005"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.703.47"
$ws.Range("E21").Value = "  +2.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.278"
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.11"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.104"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.008.93"
$ws.Range("E25").Value = "  -1.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.83"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.33"
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.868"
$ws.Range("E28").Value = "  +3.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.987"
$ws.Range("E29").Value = "  +1.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.20"
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09337"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9192"
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.297"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.336"
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.327"
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05804"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.152"
$ws.Range("E37").Value = "  +1.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02079"
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5628"
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1786"
$ws.Range("E41").Value = "  +1.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.760"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07233"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5284"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.149"
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.123"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.841"
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.81"
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.415"
$ws.Range("E50").Value = "  +4.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.003"
$ws.Range("E51").Value = "  +0.30%  "
